# Weekly update: insert a new record as row 301 (Hortaliza / Haba,
# Mercado Mayorista Lo Valledor de Santiago), pushing every existing
# row from 301 downward down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 301 - shifts rows 301:409 down to 302:410
# and copies formatting (incl. the date style on column D) from the
# row above, same as Excel's native "Insert Sheet Rows" command.
$ws.Rows("301:301").Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A301").Value2 = 6
$ws.Range("B301").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C301").Value2 = "Metropolitana"
$ws.Range("D301").Value2 = 45119
$ws.Range("E301").Value2 = 13
$ws.Range("F301").Value2 = 100112026
$ws.Range("G301").Value2 = "Haba"
$ws.Range("H301").Value2 = "Sin especificar"
$ws.Range("I301").Value2 = "Primera"
$ws.Range("J301").Value2 = 520
$ws.Range("K301").Value2 = 14000
$ws.Range("L301").Value2 = 15000
$ws.Range("M301").Value2 = 14385
$ws.Range("N301").Value2 = "`$/saco 25 kilos"
$ws.Range("O301").Value2 = "Provincia de Limarí"
$ws.Range("P301").Value2 = 575
$ws.Range("Q301").Value2 = 25
$ws.Range("R301").Value2 = "Hortaliza"
